# Add a "Spain" worksheet (copy of "Italy") with Spain-specific test data,
# matching the "Added Test data for Spain Zettler Market" commit.

$wb = $excel.ActiveWorkbook

# Italy is the template for every other "market" sheet in this workbook -
# copy it, place the copy right after it, then rename/re-point the copy.
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Market-specific content.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034"

# The Spain sheet uses taller rows for the wrapped user-story text.
$spain.Range("A3:D3").RowHeight = 28.8
$spain.Range("A4:D4").RowHeight = 28.8

# Column widths tailored to the new content (closest achievable values;
# column C keeps the width inherited from the Italy template).
$spain.Columns.Item(1).ColumnWidth = 24.3
$spain.Columns.Item(2).ColumnWidth = 35.85
$spain.Columns.Item(4).ColumnWidth = 19.65

# Selection state: Italy is no longer the active tab, selection becomes the
# full used range; Spain becomes the active tab with C4 selected.
$italy.Range("A1:D11").Select()
$spain.Activate()
$spain.Range("C4").Select()
